$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("peat")

$ws.Range("C1").Value = "zeta_diri_BC"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "-0.2"
$ws.Range("C2").NumberFormat = "General"

$ws.Range("C3").Select()
